# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de sheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 13:04:12"
$wsZhCn.Range("H2").Value = "2016-03-13 13:04:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 13:04:16"
$wsDeDe.Range("H2").Value = "2016-03-13 13:04:36"
